$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.322.91'
$ws.Range("E2").Value = '  -3.48%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.938.10'
$ws.Range("E3").Value = '  -3.19%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.97'
$ws.Range("E5").Value = '  -3.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7185'
$ws.Range("E6").Value = '  -11.66%  '
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3285'
$ws.Range("E8").Value = '  -8.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '26.74'
$ws.Range("E9").Value = '  +2.98%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06856'
$ws.Range("E10").Value = '  -2.62%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8093'
$ws.Range("E11").Value = '  -4.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07991'
$ws.Range("E12").Value = '  -1.98%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.934.85'
$ws.Range("E13").Value = '  -3.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.455'
$ws.Range("E14").Value = '  -1.47%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '95.16'
$ws.Range("E15").Value = '  -6.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.64'
$ws.Range("E16").Value = '  +4.06%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '264.52'
$ws.Range("E17").Value = '  -3.47%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '30.321.78'
$ws.Range("E18").Value = '  -3.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.870'
$ws.Range("E19").Value = '  +0.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007995'
$ws.Range("E20").Value = '  +0.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.187.67'
$ws.Range("E21").Value = '  -2.81%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  +0.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9997'
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.950'
$ws.Range("E24").Value = '  -2.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.799'
$ws.Range("E25").Value = '  -0.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.43'
$ws.Range("E26").Value = '  -2.68%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.353'
$ws.Range("E27").Value = '  +3.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1349'
$ws.Range("E28").Value = '  -12.94%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.09'
$ws.Range("E29").Value = '  -5.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.366'
$ws.Range("E30").Value = '  -0.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.558'
$ws.Range("E31").Value = '  -1.38%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.435'
$ws.Range("E32").Value = '  -5.09%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.256'
$ws.Range("E33").Value = '  -2.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05110'
$ws.Range("E34").Value = '  -2.43%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.213'
$ws.Range("E35").Value = '  -0.79%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7508'
$ws.Range("E36").Value = '  -1.94%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.742'
$ws.Range("E37").Value = '  -0.52%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01951'
$ws.Range("E38").Value = '  -3.46%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.820'
$ws.Range("E39").Value = '  -3.43%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '81.36'
$ws.Range("E40").Value = '  +2.93%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.597'
$ws.Range("E41").Value = '  -1.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4507'
$ws.Range("E42").Value = '  -5.78%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.031'
$ws.Range("E43").Value = '  -5.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8395'
$ws.Range("E45").Value = '  -2.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.62'
$ws.Range("E46").Value = '  -2.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.771'
$ws.Range("E47").Value = '  -2.60%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.374'
$ws.Range("E48").Value = '  -2.50%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.40'
$ws.Range("E49").Value = '  -1.63%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.503'
$ws.Range("E50").Value = '  +2.67%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4142'
$ws.Range("E51").Value = '  -6.15%  '
